# "Add Few More TestCases And Refactor Structure"
#
# - Adds a new row (A3) to the "boostrapAlertSuccess" sheet (2nd sheet) with
#   the text "Faile for Raport Purpoose" (new shared string).
# - Makes that sheet the active/selected sheet (was previously the last
#   sheet, "boostrapAlertInfo", that was active/selected).

$wb = $excel.ActiveWorkbook

# The sheet that receives the new row is "boostrapAlertSuccess" (2nd tab).
$ws = $wb.Worksheets.Item(2)

# Activating this sheet makes it the workbook's active tab and marks its
# window as the selected/visible tab (and correspondingly un-marks whatever
# sheet used to be selected, e.g. the last sheet).
$ws.Activate()

# New data row under the existing two rows.
$ws.Range("A3").Value = "Faile for Raport Purpoose"

# Leave the new cell selected, matching the saved selection state.
$ws.Range("A3").Select()
